$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'306.89"
$ws.Range("E2").Value = "'7.04%"

# Row 3
$ws.Range("D3").Value = "'32.34"
$ws.Range("E3").Value = "'10.21%"

# Row 4
$ws.Range("D4").Value = "'5.328"
$ws.Range("E4").Value = "'4.80%"

# Row 5
$ws.Range("D5").Value = "'0.07422"
$ws.Range("E5").Value = "'11.47%"

# Row 6
$ws.Range("D6").Value = "'7.781"
$ws.Range("E6").Value = "'6.22%"

# Row 7
$ws.Range("D7").Value = "'3.692"
$ws.Range("E7").Value = "'8.42%"

# Row 8
$ws.Range("D8").Value = "'1.585"
$ws.Range("E8").Value = "'17.77%"

# Row 9
$ws.Range("D9").Value = "'0.9133"
$ws.Range("E9").Value = "'-0.87%"

# Row 10
$ws.Range("D10").Value = "'0.01652"
$ws.Range("E10").Value = "'2,460.97%"

# Row 11
$ws.Range("D11").Value = "'0.1670"
$ws.Range("E11").Value = "'7.19%"

# Row 12
$ws.Range("D12").Value = "'0.07414"
$ws.Range("E12").Value = "'15.80%"

# Row 13
$ws.Range("D13").Value = "'0.08051"
$ws.Range("E13").Value = "'6.42%"

# Row 14
$ws.Range("D14").Value = "'0.03023"
$ws.Range("E14").Value = "'4.21%"

# Row 15
$ws.Range("D15").Value = "'0.09815"
$ws.Range("E15").Value = "'9.20%"

# Row 16
$ws.Range("D16").Value = "'0.001517"
$ws.Range("E16").Value = "'-4.29%"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04561"
$ws.Range("E17").Value = "'1.53%"

# Row 18
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006106"
$ws.Range("E18").Value = "'-2.84%"

# Row 19
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.485"
$ws.Range("E19").Value = "'0.77%"

# Row 20
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").Value = "'2.248"
$ws.Range("E20").Value = "'0.83%"

# Row 21
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").Value = "'0.3273"
$ws.Range("E21").Value = "'1.91%"

# Row 22
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1307"
$ws.Range("E22").Value = "'-0.16%"

# Row 23
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").Value = "'4.254"
$ws.Range("E23").Value = "'4.62%"

# Row 24
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").Value = "'0.1613"
$ws.Range("E24").Value = "'4.17%"

# Row 25
$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'2.31%"

# Row 26
$ws.Range("D26").Value = "'0.004534"
$ws.Range("E26").Value = "'10.12%"

# Row 27
$ws.Range("D27").Value = "'0.0001165"
$ws.Range("E27").Value = "'-6.77%"

# Row 28
$ws.Range("D28").Value = "'0.0001705"
$ws.Range("E28").Value = "'5.43%"

# Row 40
$ws.Range("D40").Value = "'0.04534"
$ws.Range("E40").Value = "'8.53%"

# Row 41
$ws.Range("D41").Value = "'0.007264"
$ws.Range("E41").Value = "'7.93%"

# Row 42
$ws.Range("D42").Value = "'0.1367"
$ws.Range("E42").Value = "'10.49%"

# Row 43
$ws.Range("D43").Value = "'0.002171"
$ws.Range("E43").Value = "'9.67%"

# Row 44
$ws.Range("D44").Value = "'0.01375"
$ws.Range("E44").Value = "'10.36%"

# Row 45
$ws.Range("D45").Value = "'0.00005937"
$ws.Range("E45").Value = "'5.76%"

# Row 46
$ws.Range("D46").Value = "'1.882"
$ws.Range("E46").Value = "'-4.34%"
